$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format so that numeric-looking
# strings (e.g. "1.008", "24.091.37", "0.00001325") are written back as
# text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.091.37"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.659.62"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").Value = "311.03"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "0.3947"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("D8").Value = "0.3927"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("D9").Value = "51.92"
$ws.Range("E9").Value = "  +4.71%  "
$ws.Range("D10").Value = "1.375"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "0.08529"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "24.26"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").Value = "7.261"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "7.996"
$ws.Range("E15").Value = "  +7.19%  "
$ws.Range("D16").Value = "0.00001325"
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("D17").Value = "1.665.50"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "95.49"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "0.07010"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "20.32"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "6.982"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "13.83"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "24.116.81"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "2.509"
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("D26").Value = "3.089"
$ws.Range("E26").Value = "  +9.66%  "
$ws.Range("D27").Value = "22.44"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "156.61"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "141.61"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "5.342"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "7.962"
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("D32").Value = "2.527"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").Value = "1.848.85"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "1.048"
$ws.Range("E34").Value = "  +10.14%  "
$ws.Range("D35").Value = "0.03042"
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("D36").Value = "0.08171"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "6.799"
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").Value = "11.10"
$ws.Range("E38").Value = "  +10.80%  "
$ws.Range("D39").Value = "0.2741"
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("D40").Value = "0.09249"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "13.72"
$ws.Range("E41").Value = "  +5.51%  "
$ws.Range("D42").Value = "0.7620"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "1.439"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "16.55"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("D45").Value = "0.7016"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "2.511"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "4.111"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "0.08349"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "135.90"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").Value = "1.415"
$ws.Range("E51").Value = "  +7.28%  "

# Restore the default "Normal" style so no stray style index is left on
# these cells (matching the original workbook formatting).
$ws.Range("D2:E51").Style = "Normal"
